# Insert a new colo row for Grenada (GND / St. George's) directly above the
# existing "AMM" (Amman, Jordan) row, pushing AMM and every row below it down
# by one. This grows the used range from A1:G310 to A1:G311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 216; everything from 216 downward shifts to 217+.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row with the Grenada colo entry.
$ws.Range("A216").Value = "GND"
$ws.Range("B216").Value = "St. George's, Grenada"
$ws.Range("C216").Value = 12.007116
$ws.Range("D216").Value = -61.7882288
$ws.Range("E216").Value = "GD"
$ws.Range("F216").Value = "South America"
$ws.Range("G216").Value = "St. George's"

# Match the formatting of the "colo" column (bold, centered, bordered) used by
# every other row -- copy it from the row right below (the shifted AMM row).
$ws.Range("A217").Copy()
$ws.Range("A216").PasteSpecial(-4122)
